# Update "想去人数" (interest count) figures on the "展览" and "全部类型" sheets
# to reflect the refreshed scrape data (gh-pages output at 456a3b4):
#   F2: 151 -> 152
#   F4: 249 -> 250
#   F5: 3901 -> 3913

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 152
    $ws.Range("F4").Value = 250
    $ws.Range("F5").Value = 3913
}
